$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "pitstop" sheet: add two new columns H (tyre_before) and I (tyre_after)
# ---------------------------------------------------------------------------
$pitstop = $wb.Worksheets.Item("pitstop")

$pitstop.Range("H1").Value = "tyre_before"
$pitstop.Range("I1").Value = "tyre_after"

$tyreData = @(
    @(2,2,3),@(3,3,3),@(4,3,1),@(5,1,2),@(6,3,2),@(7,2,3),@(8,3,2),@(9,2,3),
    @(10,3,2),@(11,2,3),@(12,3,2),@(13,2,3),@(14,3,2),@(15,2,3),@(16,2,1),@(17,1,3),
    @(18,3,2),@(19,2,3),@(20,3,2),@(21,2,3),@(22,3,2),@(23,2,3),@(24,2,3),@(25,3,3),
    @(26,2,1),@(27,1,3),@(28,3,3),@(29,3,3),@(30,3,2),@(31,2,3),@(32,3,3),@(33,3,2),
    @(34,2,3),@(35,3,2),@(36,2,2),@(37,2,3),@(38,3,3),@(39,3,1),@(40,1,2),@(41,2,2),
    @(42,2,1),@(43,1,3),@(44,3,3),@(45,3,2),@(46,2,3)
)

foreach ($row in $tyreData) {
    $r = $row[0]
    $before = $row[1]
    $after = $row[2]
    $pitstop.Cells.Item($r, 8).Value = $before
    $pitstop.Cells.Item($r, 9).Value = $after
}

$pitstop.Columns.Item(8).AutoFit()
$pitstop.Columns.Item(9).AutoFit()
$pitstop.Range("H47").Select()

# ---------------------------------------------------------------------------
# 2) Two brand-new sheets at the end of the workbook: "weather" and "altitude"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$weather = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$weather.Name = "weather"
$altitude = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $weather)
$altitude.Name = "altitude"

$weather.Range("A1").Value = "Skycondition"
$weather.Range("B1").Value = "Partly Cloudy"
$weather.Range("A2").Value = "Temperature"
$weather.Range("B2").Value = "67.2°F"
$weather.Range("A3").Value = "Humidity"
$weather.Range("B3").Value = 0.66
$weather.Range("B3").NumberFormat = "0%"
$weather.Range("A4").Value = "Wind speed"
$weather.Range("B4").Value = "14.98 mph"
$weather.Range("A5").Value = "Wind bearing"
$weather.Range("B5").Value = "164°"
$weather.Range("C5").Select()

$altitude.Range("A1").Value = "delta"
$altitude.Range("B1").Value = 43
$altitude.Range("B2").Select()
